$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") has some cells whose value is the combined
# string "System, dnasr281@gmail.com". The author re-ordered that list to
# "dnasr281@gmail.com, System". Do an exact-match replace scoped to the
# used range so only full matches are affected (cells that contain just
# "System" or just "dnasr281@gmail.com" are left untouched).

$lastRow = $ws.UsedRange.Rows.Count
$rangeG = $ws.Range("G1:G$lastRow")

$rangeG.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
